# "corrections to overlay and added overlay to main"
# Adds a new worksheet named "5" (mirroring the existing "10 second wait for
# orders" / "20" / "10" overlay sheets) with 5 days of data, positioned as the
# last sheet in the workbook.

$wb = $excel.ActiveWorkbook

# Use the "10" sheet as a template: a full Worksheet.Copy duplicates the
# sheet (styles, column widths, page setup, etc.) exactly, which is more
# faithful than rebuilding formatting cell-by-cell.
$template = $wb.Worksheets.Item("10")
$template.Copy($null, $template)

# The freshly duplicated sheet is placed right after the template and is
# named "10 (2)" (or similar); grab it and rename it to "5".
$newSheet = $wb.Worksheets.Item($template.Index + 1)
$newSheet.Name = "5"

# Overlay data for the new "5" sheet.
$data = @(
    @(1, 55.94,  36, 199,  1, 2,  0, 0, 0, 2),
    @(2, 522.28, 36, 588,  2, 5,  0, 0, 0, 5),
    @(3, 515.92, 27, 982,  3, 8,  0, 0, 0, 8),
    @(4, 530.08, 25, 1468, 4, 10, 0, 0, 0, 10),
    @(5, 576.42, 17, 1957, 4, 12, 0, 0, 0, 12)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $newSheet.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# The template had 10 data rows (rows 2-11); the new sheet only needs 5
# (rows 2-6), so clear out the leftover rows 7-11 entirely.
$newSheet.Range("A7:J11").Clear()

# Restore the original active/selected sheet so tabSelected stays put.
$wb.Worksheets.Item("10 second wait for orders").Activate()
